$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("H4").Value = '2/217'
$ws.Range("G7").Value = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = '48.5%'
$ws.Range("G15").Value = 'dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = '52.7%'
$ws.Range("S16").NumberFormat = "@"
$ws.Range("S16").Value = '60.4%'
$ws.Range("G17").Value = 'Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G18").Value = 'marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G21").Value = 'aya.saeed@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg'
$ws.Range("H21").Value = '166/216'
$ws.Range("G24").Value = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G32").Value = 'dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range("G34").Value = 'Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg'
$ws.Range("G35").Value = 'Aya_hamed@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg'
$ws.Range("G37").Value = 'System, Kerelos.zareef@med.asu.edu.eg'
$ws.Range("G45").Value = 'mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com'
$ws.Range("G48").Value = 'afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg'
$ws.Range("G51").Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G52").Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G54").Value = 'System, Kerelos.zareef@med.asu.edu.eg'
$ws.Range("G62").Value = 'mohamed.saleem@med.asu.edu.eg, System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com'
$ws.Range("G65").Value = 'afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg, Aya.alaa-eldein@med.asu.edu.eg'
$ws.Range("G68").Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G69").Value = 'naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G72").Value = 'mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg'
$ws.Range("G76").Value = 'mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg'
$ws.Range("G83").Value = 'Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range("G85").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G86").Value = 'marina_atef@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg'
$ws.Range("G98").Value = 'nourhanhosni@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, user@user.com, nourhanmohamed@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G99").Value = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G102").Value = 'Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G109").Value = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G116").Value = 'afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G119").Value = 'marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G126").Value = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G133").Value = 'afaf.abdallah@med.asu.edu.eg, enas.omran@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg'
$ws.Range("G136").Value = 'marina_atef@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, marinasorial@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg'
$ws.Range("G143").Value = 'mohamed.saleem@med.asu.edu.eg, afnan.fares@med.asu.edu.eg'
$ws.Range("G144").Value = 'mohamed.saleem@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, nourhanmamdouh@med.asu.edu.eg'
$ws.Range("G149").Value = 'user@user.com, Walaa.h.ghanima@med.asu.edu.eg'
$ws.Range("G150").Value = 'Youstina.ibrahim@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg'
$ws.Range("G151").Value = 'nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg, System, randa.rabea@med.asu.edu.eg'

$ws.Range("G117").Value = "user@user.com, Amr-Saeed@med.asu.edu.eg`n"
